$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (dates as Excel serial numbers), continuing the
# existing series through 2021-09-20 (commit: "aggiornamento fino a 20/09/2021").
$data = @(
    @(44449, 1, 8, 129.3870289503477),
    @(44450, 4, 12, 194.0805434255216),
    @(44451, 1, 12, 194.0805434255216),
    @(44452, 0, 9, 145.5604075691412),
    @(44453, 0, 7, 113.2136503315543),
    @(44454, 0, 7, 113.2136503315543),
    @(44455, 1, 7, 113.2136503315543),
    @(44456, 1, 7, 113.2136503315543),
    @(44457, 0, 3, 48.5201358563804),
    @(44458, 0, 2, 32.34675723758694),
    @(44459, 0, 2, 32.34675723758694)
)

$startRow = 375
$endRow = $startRow + $data.Length - 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

# Carry the date-column style (border/format) from the last existing
# row down onto the newly appended rows, same as Excel's "fill down".
$ws.Cells.Item($startRow - 1, 1).Copy()
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0
